$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Fix employee ID ---
# (shared string emp_fn0y5dge -> emp_emnnysju), found on "Jason Schema" sheet column B
$wsSchema.Range("B2:B6").Value = "emp_emnnysju"

# --- Fix client names (Weekly Timesheet column B, rows 2-6) ---
$wsTime.Range("B2").Value = "Prezzano"
$wsTime.Range("B3").Value = "Vincent"
$wsTime.Range("B4").Value = "Zygmunt"
$wsTime.Range("B5").Value = "Ricca"
$wsTime.Range("B6").Value = "Varricchio"

# --- Mirror client names on Jason Schema sheet column D, rows 2-6 ---
$wsSchema.Range("D2").Value = "Prezzano"
$wsSchema.Range("D3").Value = "Vincent"
$wsSchema.Range("D4").Value = "Zygmunt"
$wsSchema.Range("D5").Value = "Ricca"
$wsSchema.Range("D6").Value = "Varricchio"

# --- Simulator full-month coverage: populate Rate/Total columns ---
# Weekly Timesheet: E2:E6 (Rate) = 88, F2:F6 (Total) = 704
for ($r = 2; $r -le 6; $r++) {
    $wsTime.Cells.Item($r, 5).Value = 88
    $wsTime.Cells.Item($r, 6).Value = 704
}

# Subtotal rows on Weekly Timesheet
$wsTime.Range("F8").Value = 3520
$wsTime.Range("F11").Value = 3520
$wsTime.Range("F13").Value = 3520

# Jason Schema: F2:F6 (Rate) = 88, G2:G6 (Total) = 704
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 6).Value = 88
    $wsSchema.Cells.Item($r, 7).Value = 704
}
